# Lua console: add "run file" + "brush" features.
# - Drop the old description row, the old "LuaScriptConsoleExecute"/"运行脚本"
#   row, and the old "LuaScriptConsoleRun" row (its value changes).
# - Append four rows at the bottom: the (renamed-value) Run row, a new
#   RunFile row, a new Brush row, and the updated Description row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so row numbers of the earlier deletions stay valid.
$ws.Rows(571).Delete()   # LuaScriptConsoleRun / 执行输入        (value changes -> re-added below)
$ws.Rows(570).Delete()   # LuaScriptConsoleExecute / 运行脚本    (removed entirely)
$ws.Rows(566).Delete()   # LuaScriptConsoleDescription / (old text) (re-added below with new text)

# The sheet now has 603 data rows (566..603 reused by the former 567..606).
# Copy the formatting (styles s="4"/s="3") of the last existing row onto the
# four new rows that get appended at the end.
$ws.Range("A603:B603").Copy()
$ws.Range("A604:B607").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A604").Value2 = "LuaScriptConsoleRun"
$ws.Range("B604").Value2 = "执行"

$ws.Range("A605").Value2 = "LuaScriptConsoleRunFile"
$ws.Range("B605").Value2 = "执行文件"

$ws.Range("A606").Value2 = "LuaScriptConsoleBrush"
$ws.Range("B606").Value2 = "脚本刷"

$ws.Range("A607").Value2 = "LuaScriptConsoleDescription"
$ws.Range("B607").Value2 = "勾选“执行文件”以读取选中的Lua脚本，否则读取输入窗口。点击“执行”以执行所选的代码。点击“脚本刷”在地图的指定坐标执行脚本。脚本可能会损坏地图，请在运行前保存或执行快照函数。请参阅文档了解可用函数。"

$ws.Rows(607).RowHeight = 70

# Match the final selection recorded in the workbook view.
$ws.Range("B602").Select()
